# Auto-generated script applying scheduled market-price-refresh updates
# to the Lich_Profits workbook, per commit 'chore: update Sheets via scheduled runner'.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 1060  # H9: 1075.25 -> 1060
$ws.Cells.Item(9, 10).Value = 999  # J9: 0 -> 999
$ws.Cells.Item(9, 12).Value = 999  # L9: 0 -> 999
$ws.Cells.Item(9, 14).Value = -1337  # N9: None -> -1337
$ws.Cells.Item(15, 8).Value = 2012.7354  # H15: 1870.619 -> 2012.7354
$ws.Cells.Item(15, 9).Value = 2012.7354  # I15: 1870.619 -> 2012.7354
$ws.Cells.Item(15, 11).Value = 6038.206200000001  # K15: 5611.857 -> 6038.206200000001
$ws.Cells.Item(15, 13).Value = -5869.206200000001  # M15: -5442.857 -> -5869.206200000001
$ws.Cells.Item(40, 8).Value = 2880.0908  # H40: 3826.3572 -> 2880.0908
$ws.Cells.Item(40, 9).Value = 3198.2222  # I40: 5197.4443 -> 3198.2222
$ws.Cells.Item(40, 10).Value = 1448.5  # J40: 1358.4 -> 1448.5
$ws.Cells.Item(40, 11).Value = 3198.2222  # K40: 5197.4443 -> 3198.2222
$ws.Cells.Item(40, 12).Value = 1448.5  # L40: 1358.4 -> 1448.5
$ws.Cells.Item(40, 13).Value = -3023.2222  # M40: -5022.4443 -> -3023.2222
$ws.Cells.Item(40, 14).Value = -1798.5  # N40: -1708.4 -> -1798.5
$ws.Cells.Item(42, 8).Value = 257.07693  # H42: 258.69232 -> 257.07693
$ws.Cells.Item(42, 9).Value = 260.33334  # I42: 241.3 -> 260.33334
$ws.Cells.Item(42, 10).Value = 249.75  # J42: 316.66666 -> 249.75
$ws.Cells.Item(42, 11).Value = 781.0000200000001  # K42: 723.9000000000001 -> 781.0000200000001
$ws.Cells.Item(42, 12).Value = 749.25  # L42: 949.9999799999999 -> 749.25
$ws.Cells.Item(42, 13).Value = -551.0000200000001  # M42: -493.9000000000001 -> -551.0000200000001
$ws.Cells.Item(42, 14).Value = -1209.25  # N42: -1409.99998 -> -1209.25
$ws.Cells.Item(43, 8).Value = 3596.3333  # H43: 4232.6665 -> 3596.3333
$ws.Cells.Item(43, 9).Value = 2994.5  # I43: 2999 -> 2994.5
$ws.Cells.Item(43, 10).Value = 4800  # J43: 4849.5 -> 4800
$ws.Cells.Item(43, 11).Value = 2994.5  # K43: 2999 -> 2994.5
$ws.Cells.Item(43, 12).Value = 4800  # L43: 4849.5 -> 4800
$ws.Cells.Item(43, 13).Value = -2925.5  # M43: -2930 -> -2925.5
$ws.Cells.Item(43, 14).Value = -4938  # N43: -4987.5 -> -4938
$ws.Cells.Item(94, 8).Value = 1500.8  # H94: 1367.1666 -> 1500.8
$ws.Cells.Item(94, 9).Value = 1500.8  # I94: 1367.1666 -> 1500.8
$ws.Cells.Item(94, 11).Value = 1500.8  # K94: 1367.1666 -> 1500.8
$ws.Cells.Item(94, 13).Value = -1049.8  # M94: -916.1666 -> -1049.8
$ws.Cells.Item(111, 8).Value = 3119.2144  # H111: 3289.923 -> 3119.2144
$ws.Cells.Item(111, 10).Value = 3857.111  # J111: 4226.75 -> 3857.111
$ws.Cells.Item(111, 12).Value = 11571.333  # L111: 12680.25 -> 11571.333
$ws.Cells.Item(111, 14).Value = -17705.333  # N111: -18814.25 -> -17705.333
$ws.Cells.Item(112, 8).Value = 5439708.5  # H112: 5892910 -> 5439708.5
$ws.Cells.Item(112, 9).Value = 1471.8  # I112: 1517.5 -> 1471.8
$ws.Cells.Item(112, 11).Value = 4415.4  # K112: 4552.5 -> 4415.4
$ws.Cells.Item(112, 13).Value = -3307.4  # M112: -3444.5 -> -3307.4
$ws.Cells.Item(113, 8).Value = 6809.56  # H113: 6922.5835 -> 6809.56
$ws.Cells.Item(113, 9).Value = 7618.8887  # I113: 7826.0586 -> 7618.8887
$ws.Cells.Item(113, 11).Value = 7618.8887  # K113: 7826.0586 -> 7618.8887
$ws.Cells.Item(113, 13).Value = -4364.8887  # M113: -4572.0586 -> -4364.8887
$ws.Cells.Item(115, 8).Value = 67348936  # H115: 63607364 -> 67348936
$ws.Cells.Item(115, 9).Value = 76327464  # I115: 71557030 -> 76327464
$ws.Cells.Item(115, 11).Value = 228982392  # K115: 214671090 -> 228982392
$ws.Cells.Item(115, 13).Value = -228980825  # M115: -214669523 -> -228980825
$ws.Cells.Item(127, 8).Value = 5785.6924  # H127: 6465.2144 -> 5785.6924
$ws.Cells.Item(127, 9).Value = 1704  # I127: 1920.6666 -> 1704
$ws.Cells.Item(127, 10).Value = 9284.286  # J127: 9873.625 -> 9284.286
$ws.Cells.Item(127, 11).Value = 5112  # K127: 5761.9998 -> 5112
$ws.Cells.Item(127, 12).Value = 27852.858  # L127: 29620.875 -> 27852.858
$ws.Cells.Item(127, 13).Value = -152  # M127: -801.9997999999996 -> -152
$ws.Cells.Item(127, 14).Value = -37772.858  # N127: -39540.875 -> -37772.858
$ws.Cells.Item(132, 8).Value = 1357.9608  # H132: 1324.5636 -> 1357.9608
$ws.Cells.Item(132, 9).Value = 1099.3414  # I132: 1077.8334 -> 1099.3414
$ws.Cells.Item(132, 10).Value = 2418.3  # J132: 2121.6924 -> 2418.3
$ws.Cells.Item(132, 11).Value = 3298.0242  # K132: 3233.5002 -> 3298.0242
$ws.Cells.Item(132, 12).Value = 7254.900000000001  # L132: 6365.0772 -> 7254.900000000001
$ws.Cells.Item(132, 13).Value = -768.0241999999998  # M132: -703.5001999999999 -> -768.0241999999998
$ws.Cells.Item(132, 14).Value = -12314.9  # N132: -11425.0772 -> -12314.9
$ws.Cells.Item(135, 8).Value = 6180.75  # H135: 4225.857 -> 6180.75
$ws.Cells.Item(135, 9).Value = 4032  # I135: 2222.5 -> 4032
$ws.Cells.Item(135, 11).Value = 36288  # K135: 20002.5 -> 36288
$ws.Cells.Item(135, 13).Value = -33753  # M135: -17467.5 -> -33753

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(101, 8).Value = 100598.8  # H101: 100599.75 -> 100598.8
$ws.Cells.Item(101, 10).Value = 100598.8  # J101: 100599.75 -> 100598.8
$ws.Cells.Item(101, 12).Value = 100598.8  # L101: 100599.75 -> 100598.8
$ws.Cells.Item(101, 14).Value = -107088.8  # N101: -107089.75 -> -107088.8
$ws.Cells.Item(132, 8).Value = 3223.0334  # H132: 2701.4595 -> 3223.0334
$ws.Cells.Item(132, 9).Value = 4111.7144  # I132: 3200.3215 -> 4111.7144
$ws.Cells.Item(132, 11).Value = 12335.1432  # K132: 9600.9645 -> 12335.1432
$ws.Cells.Item(132, 13).Value = -9805.143199999999  # M132: -7070.9645 -> -9805.143199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3377.0715  # H99: 3309.848 -> 3377.0715
$ws.Cells.Item(99, 9).Value = 2976.4443  # I99: 3028.5186 -> 2976.4443
$ws.Cells.Item(99, 10).Value = 4098.2  # J99: 3709.6316 -> 4098.2
$ws.Cells.Item(99, 11).Value = 2976.4443  # K99: 3028.5186 -> 2976.4443
$ws.Cells.Item(99, 12).Value = 4098.2  # L99: 3709.6316 -> 4098.2
$ws.Cells.Item(99, 13).Value = -1478.4443  # M99: -1530.5186 -> -1478.4443
$ws.Cells.Item(99, 14).Value = -7094.2  # N99: -6705.631600000001 -> -7094.2
$ws.Cells.Item(134, 8).Value = 2780.74  # H134: 2545.3818 -> 2780.74
$ws.Cells.Item(134, 9).Value = 2527.8096  # I134: 2279.2979 -> 2527.8096
$ws.Cells.Item(134, 11).Value = 7583.4288  # K134: 6837.893700000001 -> 7583.4288
$ws.Cells.Item(134, 13).Value = -5048.4288  # M134: -4302.893700000001 -> -5048.4288

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(43, 8).Value = 21326  # H43: 21091.2 -> 21326
$ws.Cells.Item(43, 10).Value = 21326  # J43: 21091.2 -> 21326
$ws.Cells.Item(43, 12).Value = 21326  # L43: 21091.2 -> 21326
$ws.Cells.Item(43, 14).Value = -21694  # N43: -21459.2 -> -21694
$ws.Cells.Item(55, 8).Value = 7326.3335  # H55: 7326.6665 -> 7326.3335
$ws.Cells.Item(55, 9).Value = 6499  # I55: 6499.5 -> 6499
$ws.Cells.Item(55, 11).Value = 6499  # K55: 6499.5 -> 6499
$ws.Cells.Item(55, 13).Value = -6184  # M55: -6184.5 -> -6184
$ws.Cells.Item(101, 8).Value = 21326  # H101: 21091.2 -> 21326
$ws.Cells.Item(101, 10).Value = 21326  # J101: 21091.2 -> 21326
$ws.Cells.Item(101, 12).Value = 21326  # L101: 21091.2 -> 21326
$ws.Cells.Item(101, 14).Value = -27816  # N101: -27581.2 -> -27816
$ws.Cells.Item(132, 8).Value = 3460.125  # H132: 3496.8572 -> 3460.125
$ws.Cells.Item(132, 9).Value = 3597.2856  # I132: 3496.8572 -> 3597.2856
$ws.Cells.Item(132, 10).Value = 2500  # J132: 0 -> 2500
$ws.Cells.Item(132, 11).Value = 10791.8568  # K132: 10490.5716 -> 10791.8568
$ws.Cells.Item(132, 12).Value = 7500  # L132: 0 -> 7500
$ws.Cells.Item(132, 13).Value = -8261.856800000001  # M132: -7960.571599999999 -> -8261.856800000001
$ws.Cells.Item(132, 14).Value = -12560  # N132: None -> -12560
$ws.Cells.Item(134, 8).Value = 3780.2222  # H134: 3787.4443 -> 3780.2222
$ws.Cells.Item(134, 9).Value = 4310.4443  # I134: 4319.472 -> 4310.4443
$ws.Cells.Item(134, 11).Value = 12931.3329  # K134: 12958.416 -> 12931.3329
$ws.Cells.Item(134, 13).Value = -10396.3329  # M134: -10423.416 -> -10396.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 2999.5  # H11: 4100 -> 2999.5
$ws.Cells.Item(11, 9).Value = 2999.5  # I11: 4100 -> 2999.5
$ws.Cells.Item(11, 11).Value = 8998.5  # K11: 12300 -> 8998.5
$ws.Cells.Item(11, 13).Value = -8858.5  # M11: -12160 -> -8858.5
$ws.Cells.Item(32, 8).Value = 30666.334  # H32: 31425 -> 30666.334
$ws.Cells.Item(32, 10).Value = 30666.334  # J32: 31425 -> 30666.334
$ws.Cells.Item(32, 12).Value = 91999.00199999999  # L32: 94275 -> 91999.00199999999
$ws.Cells.Item(32, 14).Value = -92565.00199999999  # N32: -94841 -> -92565.00199999999
$ws.Cells.Item(70, 8).Value = 1499.5  # H70: 1571 -> 1499.5
$ws.Cells.Item(73, 8).Value = 1499.5  # H73: 1571 -> 1499.5
$ws.Cells.Item(113, 8).Value = 0  # H113: 2550 -> 0
$ws.Cells.Item(113, 9).Value = 0  # I113: 2550 -> 0
$ws.Cells.Item(113, 11).Value = 0  # K113: 7650 -> 0
$ws.Cells.Item(113, 13).ClearContents()  # M113: -5480 -> (cleared)
$ws.Cells.Item(122, 8).Value = 1787.125  # H122: 1883.1666 -> 1787.125
$ws.Cells.Item(122, 10).Value = 1732.1666  # J122: 1848.75 -> 1732.1666
$ws.Cells.Item(122, 12).Value = 15589.4994  # L122: 16638.75 -> 15589.4994
$ws.Cells.Item(122, 14).Value = -20489.4994  # N122: -21538.75 -> -20489.4994
$ws.Cells.Item(131, 8).Value = 6025541  # H131: 5557001 -> 6025541
$ws.Cells.Item(131, 10).Value = 1462.1464  # J131: 1461.663 -> 1462.1464
$ws.Cells.Item(131, 12).Value = 4386.439200000001  # L131: 4384.989 -> 4386.439200000001
$ws.Cells.Item(131, 14).Value = -14466.4392  # N131: -14464.989 -> -14466.4392

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 125  # H2: 122 -> 125
$ws.Cells.Item(2, 10).Value = 150  # J2: 136.66667 -> 150
$ws.Cells.Item(2, 12).Value = 150  # L2: 136.66667 -> 150
$ws.Cells.Item(2, 14).Value = -376  # N2: -362.66667 -> -376
$ws.Cells.Item(3, 8).Value = 20034.625  # H3: 14847 -> 20034.625
$ws.Cells.Item(3, 9).Value = 11431  # I3: 6958.4 -> 11431
$ws.Cells.Item(3, 10).Value = 25196.8  # J3: 27994.666 -> 25196.8
$ws.Cells.Item(3, 11).Value = 11431  # K3: 6958.4 -> 11431
$ws.Cells.Item(3, 12).Value = 25196.8  # L3: 27994.666 -> 25196.8
$ws.Cells.Item(3, 13).Value = -11315  # M3: -6842.4 -> -11315
$ws.Cells.Item(3, 14).Value = -25428.8  # N3: -28226.666 -> -25428.8
$ws.Cells.Item(6, 8).Value = 4500  # H6: 4250 -> 4500
$ws.Cells.Item(6, 10).Value = 4500  # J6: 4250 -> 4500
$ws.Cells.Item(6, 12).Value = 4500  # L6: 4250 -> 4500
$ws.Cells.Item(6, 14).Value = -4726  # N6: -4476 -> -4726
$ws.Cells.Item(9, 8).Value = 4998.6665  # H9: 3749.25 -> 4998.6665
$ws.Cells.Item(9, 9).Value = 4999  # I9: 3332.6667 -> 4999
$ws.Cells.Item(9, 10).Value = 4998.5  # J9: 4999 -> 4998.5
$ws.Cells.Item(9, 11).Value = 4999  # K9: 3332.6667 -> 4999
$ws.Cells.Item(9, 12).Value = 4998.5  # L9: 4999 -> 4998.5
$ws.Cells.Item(9, 13).Value = -4829  # M9: -3162.6667 -> -4829
$ws.Cells.Item(9, 14).Value = -5338.5  # N9: -5339 -> -5338.5
$ws.Cells.Item(10, 8).Value = 9665.666999999999  # H10: 10888 -> 9665.666999999999
$ws.Cells.Item(10, 9).Value = 9499.5  # I10: 11184.667 -> 9499.5
$ws.Cells.Item(10, 11).Value = 9499.5  # K10: 11184.667 -> 9499.5
$ws.Cells.Item(10, 13).Value = -9330.5  # M10: -11015.667 -> -9330.5
$ws.Cells.Item(14, 8).Value = 506000  # H14: 8333 -> 506000
$ws.Cells.Item(14, 9).Value = 0  # I14: 1004 -> 0
$ws.Cells.Item(14, 10).Value = 506000  # J14: 11997.5 -> 506000
$ws.Cells.Item(14, 11).Value = 0  # K14: 1004 -> 0
$ws.Cells.Item(14, 12).Value = 506000  # L14: 11997.5 -> 506000
$ws.Cells.Item(14, 13).ClearContents()  # M14: -836 -> (cleared)
$ws.Cells.Item(14, 14).Value = -506336  # N14: -12333.5 -> -506336
$ws.Cells.Item(16, 8).Value = 4500  # H16: 4250 -> 4500
$ws.Cells.Item(16, 10).Value = 4500  # J16: 4250 -> 4500
$ws.Cells.Item(16, 12).Value = 4500  # L16: 4250 -> 4500
$ws.Cells.Item(16, 14).Value = -5000  # N16: -4750 -> -5000
$ws.Cells.Item(33, 8).Value = 17450  # H33: 11966.667 -> 17450
$ws.Cells.Item(33, 9).Value = 16000  # I33: 8500 -> 16000
$ws.Cells.Item(33, 11).Value = 16000  # K33: 8500 -> 16000
$ws.Cells.Item(33, 13).Value = -15748  # M33: -8248 -> -15748
$ws.Cells.Item(40, 8).Value = 10000  # H40: 0 -> 10000
$ws.Cells.Item(40, 10).Value = 10000  # J40: 0 -> 10000
$ws.Cells.Item(40, 12).Value = 10000  # L40: 0 -> 10000
$ws.Cells.Item(40, 14).Value = -10302  # N40: None -> -10302
$ws.Cells.Item(52, 8).Value = 47747  # H52: 47829.668 -> 47747
$ws.Cells.Item(52, 10).Value = 47747  # J52: 47829.668 -> 47747
$ws.Cells.Item(52, 12).Value = 47747  # L52: 47829.668 -> 47747
$ws.Cells.Item(52, 14).Value = -48265  # N52: -48347.668 -> -48265
$ws.Cells.Item(58, 8).Value = 35615.918  # H58: 27724.334 -> 35615.918
$ws.Cells.Item(58, 9).Value = 34655.75  # I58: 17999.666 -> 34655.75
$ws.Cells.Item(58, 10).Value = 37536.25  # J58: 37449 -> 37536.25
$ws.Cells.Item(58, 11).Value = 34655.75  # K58: 17999.666 -> 34655.75
$ws.Cells.Item(58, 12).Value = 37536.25  # L58: 37449 -> 37536.25
$ws.Cells.Item(58, 13).Value = -34378.75  # M58: -17722.666 -> -34378.75
$ws.Cells.Item(58, 14).Value = -38090.25  # N58: -38003 -> -38090.25
$ws.Cells.Item(101, 8).Value = 25000  # H101: 0 -> 25000
$ws.Cells.Item(101, 10).Value = 25000  # J101: 0 -> 25000
$ws.Cells.Item(101, 12).Value = 25000  # L101: 0 -> 25000
$ws.Cells.Item(101, 14).Value = -31490  # N101: None -> -31490
$ws.Cells.Item(136, 8).Value = 8223.883  # H136: 4226.533 -> 8223.883
$ws.Cells.Item(136, 10).Value = 8223.883  # J136: 4226.533 -> 8223.883
$ws.Cells.Item(136, 12).Value = 24671.649  # L136: 12679.599 -> 24671.649
$ws.Cells.Item(136, 14).Value = -29771.649  # N136: -17779.599 -> -29771.649

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1360.1666  # H22: 1287.1428 -> 1360.1666
$ws.Cells.Item(22, 9).Value = 1144.4  # I22: 1060 -> 1144.4
$ws.Cells.Item(22, 11).Value = 1144.4  # K22: 1060 -> 1144.4
$ws.Cells.Item(22, 13).Value = -849.4000000000001  # M22: -765 -> -849.4000000000001
$ws.Cells.Item(27, 8).Value = 1360.1666  # H27: 1287.1428 -> 1360.1666
$ws.Cells.Item(27, 9).Value = 1144.4  # I27: 1060 -> 1144.4
$ws.Cells.Item(27, 11).Value = 1144.4  # K27: 1060 -> 1144.4
$ws.Cells.Item(27, 13).Value = -1037.4  # M27: -953 -> -1037.4
$ws.Cells.Item(98, 8).Value = 57999.25  # H98: 59249.75 -> 57999.25
$ws.Cells.Item(98, 10).Value = 57999.25  # J98: 59249.75 -> 57999.25
$ws.Cells.Item(98, 12).Value = 57999.25  # L98: 59249.75 -> 57999.25
$ws.Cells.Item(98, 14).Value = -63989.25  # N98: -65239.75 -> -63989.25
$ws.Cells.Item(104, 8).Value = 46727.5  # H104: 49473.2 -> 46727.5
$ws.Cells.Item(104, 10).Value = 46727.5  # J104: 49473.2 -> 46727.5
$ws.Cells.Item(104, 12).Value = 46727.5  # L104: 49473.2 -> 46727.5
$ws.Cells.Item(104, 14).Value = -53715.5  # N104: -56461.2 -> -53715.5
$ws.Cells.Item(106, 8).Value = 26924.666  # H106: 28398.4 -> 26924.666
$ws.Cells.Item(106, 10).Value = 26924.666  # J106: 28398.4 -> 26924.666
$ws.Cells.Item(106, 12).Value = 26924.666  # L106: 28398.4 -> 26924.666
$ws.Cells.Item(106, 14).Value = -29448.666  # N106: -30922.4 -> -29448.666
$ws.Cells.Item(122, 8).Value = 6169.6  # H122: 6579.1763 -> 6169.6
$ws.Cells.Item(122, 9).Value = 5907.154  # I122: 6764.1113 -> 5907.154
$ws.Cells.Item(122, 10).Value = 6657  # J122: 6371.125 -> 6657
$ws.Cells.Item(122, 11).Value = 17721.462  # K122: 20292.3339 -> 17721.462
$ws.Cells.Item(122, 12).Value = 19971  # L122: 19113.375 -> 19971
$ws.Cells.Item(122, 13).Value = -15271.462  # M122: -17842.3339 -> -15271.462
$ws.Cells.Item(122, 14).Value = -24871  # N122: -24013.375 -> -24871

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 0  # H18: 3200 -> 0
$ws.Cells.Item(18, 10).Value = 0  # J18: 3200 -> 0
$ws.Cells.Item(18, 12).Value = 0  # L18: 3200 -> 0
$ws.Cells.Item(18, 14).ClearContents()  # N18: -3546 -> (cleared)
$ws.Cells.Item(95, 8).Value = 24249.5  # H95: 24999.4 -> 24249.5
$ws.Cells.Item(95, 10).Value = 24249.5  # J95: 24999.4 -> 24249.5
$ws.Cells.Item(95, 12).Value = 24249.5  # L95: 24999.4 -> 24249.5
$ws.Cells.Item(95, 14).Value = -29741.5  # N95: -30491.4 -> -29741.5
$ws.Cells.Item(104, 8).Value = 57223.5  # H104: 62299 -> 57223.5
$ws.Cells.Item(104, 10).Value = 57223.5  # J104: 62299 -> 57223.5
$ws.Cells.Item(104, 12).Value = 57223.5  # L104: 62299 -> 57223.5
$ws.Cells.Item(104, 14).Value = -64211.5  # N104: -69287 -> -64211.5
$ws.Cells.Item(108, 8).Value = 45000  # H108: 0 -> 45000
$ws.Cells.Item(108, 10).Value = 45000  # J108: 0 -> 45000
$ws.Cells.Item(108, 12).Value = 45000  # L108: 0 -> 45000
$ws.Cells.Item(108, 14).Value = -52680  # N108: None -> -52680

Write-Output "Applied 238 cell updates across 8 sheets."
